$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Help"
# ---------------------------------------------------------------------------
$help = $wb.Worksheets.Item("Help")

# A2: turn the plain instructional note into rich text - append a bold,
# 12pt warning about adding field names.
$introText = "This file will be used to upload data to the NelsonDB. "
$warnText  = "FIELD NAMES CAN BE ADDED, BUT SHOULD BE DONE SPARINGLY"
$a2 = $help.Range("A2")
$a2.Value() = $introText + $warnText
$warnChars = $a2.Characters($introText.Length + 1, $warnText.Length)
$warnChars.Font.Bold = $true
$warnChars.Font.Size = 12

# Rows 11-13: Microbe* fields used to live under the shared "lab_obsother"
# table - split them out into their own "lab_obsmicrobe" table.
$help.Range("B11:B13").Value() = "lab_obsmicrobe"
$help.Range("C13").Value() = "comments"

# Rows 15-25: the remaining Source * fields move to a new "lab_obstracker"
# table, and their example field-assignment expressions are rewritten to
# match real DB field names instead of the old "source_*" placeholders.
$help.Range("B15:B25").Value() = "lab_obstracker"
$help.Range("C15").Value() = "stock_id = Stock(seed_id).id"
$help.Range("C16").Value() = "isolate_id = Isolate(isolate_id).id"
$help.Range("C17").Value() = "obs_row_id = ObsRow(row_id).id"
$help.Range("C18").Value() = "obs_plant_id = ObsPlant(plant_id).id"
$help.Range("C19").Value() = "obs_well_id = ObsWellr(well_id).id"
$help.Range("C20").Value() = "obs_microbe_id = ObsMicrobe(microbe_id).id"
$help.Range("C21").Value() = "obs_culture_id = ObsCulture(culture_id).id"
$help.Range("C22").Value() = "obs_tissue_id = ObsTissue(tissue_id).id"
$help.Range("C23").Value() = "obs_sample_id = ObsSample(sample_id).id"
$help.Range("C24").Value() = "obs-plate_id = ObsPlate(plate_id).id"
$help.Range("C25").Value() = "obs_dna_id = ObsDNA(dna_id).id"

$help.Range("A7").Select()

# ---------------------------------------------------------------------------
# Sheet "Microbe"
# ---------------------------------------------------------------------------
$microbe = $wb.Worksheets.Item("Microbe")

# Copy the rich-text cell across so both sheets share the same string
# instead of creating a duplicate entry.
$a2.Copy()
$microbe.Range("A2").PasteSpecial()

$microbe.Range("A2").Select()
